$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").ClearFormats()
$ws.Range("D1").Value = "RF-DFS-6G"

$values = @(-63,-59,-57,-49,-47,-45,-41,-34,-28,-22,-14,-12,-11,-10,-9,-9,-10,-9,-7,-5,-3,-2,-2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Range("D1:D24").Select()
